# Update the cryptocurrency price/volume table with the latest scrape.
# Column D (Price) values are stored as literal text in the source sheet
# (e.g. "63.782.99", "1.00", "0.450") - force Text format before writing so
# the automatic number-inference does not silently reformat / truncate them
# (e.g. turn "1.00" into 1, or "588.23" into a float cell), then restore
# the default "Normal" style so no stray style index is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.782.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.081.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.99%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.080.32"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.93%  "
$ws.Range("E10").Value = "  -4.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.450"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.119"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.588.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.687.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.078.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "470.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.705"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("E25").Value = "  -6.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "80.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.20%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("E32").Value = "  -6.13%  "
$ws.Range("E33").Value = "  -8.04%  "
$ws.Range("E34").Value = "  -4.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0828"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.09%  "
$ws.Range("E36").Value = "  -3.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.96%  "
$ws.Range("E38").Value = "  -2.73%  "
$ws.Range("E39").Value = "  -5.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "432.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.16%  "
$ws.Range("E43").Value = "  -4.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.112"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.00%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0361"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.813.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.37%  "
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.82%  "
